$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the event dates between the Football row (C2) and the Cricket row (C3)
$ws.Range("C2").Value = "2023-11-02T00:00:00"
$ws.Range("C3").Value = "2023-11-05T00:00:00"

# Update the selected cell to C3
$ws.Range("C3").Select()
